$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 93.62780766666667
$ws.Range("H2").Value = 280.883423
$ws.Range("I2").Value = 0.3228593149748609
$ws.Range("J2").Value = 0.3228593149748609
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 16.034937
$ws.Range("N2").Value = 48.10481100000001
$ws.Range("O2").Value = 0.1585295454080115
$ws.Range("P2").Value = 0.1585295454080115
$ws.Range("Q2").Value = 1501.315997383117
$ws.Range("R2").Value = 13511.84397644805
$ws.Range("S2").Value = 0.05118274043370669
$ws.Range("T2").Value = 0.05118274043370669

$ws.Range("G3").Value = 93.62780766666667
$ws.Range("H3").Value = 280.883423
$ws.Range("I3").Value = 0.3228593149748609
$ws.Range("J3").Value = 0.3228593149748609
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 71.80093133333332
$ws.Range("N3").Value = 215.402794
$ws.Range("O3").Value = 0.7098605379082673
$ws.Range("P3").Value = 0.7098605379082674
$ws.Range("Q3").Value = 6722.563789164873
$ws.Range("R3").Value = 60503.07410248385
$ws.Range("S3").Value = 0.2291850869967494
$ws.Range("T3").Value = 0.2291850869967495

$ws.Range("G4").Value = 93.62780766666667
$ws.Range("H4").Value = 280.883423
$ws.Range("I4").Value = 0.3228593149748609
$ws.Range("J4").Value = 0.3228593149748609
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 13.312072
$ws.Range("N4").Value = 39.936216
$ws.Range("O4").Value = 0.1316099166837212
$ws.Range("P4").Value = 0.1316099166837212
$ws.Range("Q4").Value = 1246.380116860819
$ws.Range("R4").Value = 11217.42105174737
$ws.Range("S4").Value = 0.04249148754440473
$ws.Range("T4").Value = 0.04249148754440473

$ws.Range("G5").Value = 66.39541
$ws.Range("H5").Value = 199.18623
$ws.Range("I5").Value = 0.228953097635189
$ws.Range("J5").Value = 0.228953097635189
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 16.034937
$ws.Range("N5").Value = 48.10481100000001
$ws.Range("O5").Value = 0.1585295454080115
$ws.Range("P5").Value = 0.1585295454080115
$ws.Range("Q5").Value = 1064.64621643917
$ws.Range("R5").Value = 9581.815947952531
$ws.Range("S5").Value = 0.03629583048786258
$ws.Range("T5").Value = 0.03629583048786258

$ws.Range("G6").Value = 66.39541
$ws.Range("H6").Value = 199.18623
$ws.Range("I6").Value = 0.228953097635189
$ws.Range("J6").Value = 0.228953097635189
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 71.80093133333332
$ws.Range("N6").Value = 215.402794
$ws.Range("O6").Value = 0.7098605379082673
$ws.Range("P6").Value = 0.7098605379082674
$ws.Range("Q6").Value = 4767.252274258512
$ws.Range("R6").Value = 42905.27046832661
$ws.Range("S6").Value = 0.1625247690430793
$ws.Range("T6").Value = 0.1625247690430793

$ws.Range("G7").Value = 66.39541
$ws.Range("H7").Value = 199.18623
$ws.Range("I7").Value = 0.228953097635189
$ws.Range("J7").Value = 0.228953097635189
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.312072
$ws.Range("N7").Value = 39.936216
$ws.Range("O7").Value = 0.1316099166837212
$ws.Range("P7").Value = 0.1316099166837212
$ws.Range("Q7").Value = 883.86047838952
$ws.Range("R7").Value = 7954.74430550568
$ws.Range("S7").Value = 0.0301324981042471
$ws.Range("T7").Value = 0.0301324981042471

$ws.Range("G8").Value = 129.9724656666667
$ws.Range("H8").Value = 389.917397
$ws.Range("I8").Value = 0.4481875873899502
$ws.Range("J8").Value = 0.4481875873899502
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 16.034937
$ws.Range("N8").Value = 48.10481100000001
$ws.Range("O8").Value = 0.1585295454080115
$ws.Range("P8").Value = 0.1585295454080115
$ws.Range("Q8").Value = 2084.100298699663
$ws.Range("R8").Value = 18756.90268829697
$ws.Range("S8").Value = 0.07105097448644224
$ws.Range("T8").Value = 0.07105097448644224

$ws.Range("G9").Value = 129.9724656666667
$ws.Range("H9").Value = 389.917397
$ws.Range("I9").Value = 0.4481875873899502
$ws.Range("J9").Value = 0.4481875873899502
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 71.80093133333332
$ws.Range("N9").Value = 215.402794
$ws.Range("O9").Value = 0.7098605379082673
$ws.Range("P9").Value = 0.7098605379082674
$ws.Range("Q9").Value = 9332.144082556357
$ws.Range("R9").Value = 83989.2967430072
$ws.Range("S9").Value = 0.3181506818684386
$ws.Range("T9").Value = 0.3181506818684387

$ws.Range("G10").Value = 129.9724656666667
$ws.Range("H10").Value = 389.917397
$ws.Range("I10").Value = 0.4481875873899502
$ws.Range("J10").Value = 0.4481875873899502
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 13.312072
$ws.Range("N10").Value = 39.936216
$ws.Range("O10").Value = 0.1316099166837212
$ws.Range("P10").Value = 0.1316099166837212
$ws.Range("Q10").Value = 1730.202820972195
$ws.Range("R10").Value = 15571.82538874975
$ws.Range("S10").Value = 0.05898593103506936
$ws.Range("T10").Value = 0.05898593103506936
